$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codigos_Despesas")

# Update existing cell C38: PAGAMENTO -> SALÁRIO
$ws.Range("C38").Value = "SALÁRIO"

# Append new rows 204-207
$ws.Range("A204").Value = 204
$ws.Range("B204").Value = "ADIANTAMENTO DE SALÁRIO"
$ws.Range("C204").Value = "SALÁRIO"

$ws.Range("A205").Value = 205
$ws.Range("B205").Value = "RESCISÃO TRABALHISTA"
$ws.Range("C205").Value = "RESCISÃO TRABALHISTA"

$ws.Range("A206").Value = 206
$ws.Range("B206").Value = "OPERADOR DE MAQUINAS"
$ws.Range("C206").Value = "PRESTADOR DE SERVIÇO"

$ws.Range("A207").Value = 207
$ws.Range("B207").Value = "CENOGRAFO"
$ws.Range("C207").Value = "PRESTADOR DE SERVIÇO"

# Update the defined name to reflect new range
$wb.Names.Item("Codigos_Despesas").RefersTo = "='Codigos_Despesas'!`$A`$1:`$D`$207"
